$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 49 (id_jogo 48) ----
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "Cássio, William, Fabricio Bruno, Lucas Villalba, Kaiki, Lucas Romero, Lucas Silva, Christian, Matheus Henrique, Matheus Pereira, Kaio Jorge"
$ws.Range("G49").Value = "80, Christian, Luis Sinisterra"
$ws.Range("H49").Value = "87, Kaio Jorge, Gabigol"
$ws.Range("I49").Value = "87, Matheus Pereira, Kauã Moraes"
$ws.Range("J49").Value = "87, Matheus Henrique, Eduardo"

# ---- Row 50 (id_jogo 49) ----
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "Cássio, Kauã Moraes, Jonathan Jesus, Lucas Villalba, Kaiki, Lucas Silva, Matheus Henrique, Eduardo, Matheus Pereira, Keny Arroyo, Gabigol"
$ws.Range("C50").Value = "Gabigol"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "Matheus Pereira"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "59"

$ws.Range("F50").Value = 41
$ws.Range("G50").Value = "16, Matheus Henrique, Christian"
$ws.Range("H50").Value = "55, Jonathan Jesus, Luis Sinisterra"
$ws.Range("I50").Value = "55, Lucas Silva, João Marcelo"
$ws.Range("J50").Value = "79, Kauã Moraes, Marquinhos"
$ws.Range("K50").Value = "79, Christian, Ryan Guilherme"

# ---- Update the view/selection to match the new active cell ----
$null = $ws.Range("K50").Select()
